$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- I. Status Report: row 7 Status column flips from "In Progress" to "Completed" ---
$ws.Range("C7").Value = "Completed"

# --- III. Project Issues: row 18 was an empty placeholder row, now populated ---
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Design Document (Part II, III)"
$ws.Range("C18").Value = "In Progress"

# --- View state: scroll the frozen worksheet down a bit and land the cursor on C18 ---
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("C18").Select() | Out-Null
